$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("Duration") gets a new evidence-table data point inserted at column C
# (value 3.74155, shown with 4 decimal places). The values that used to live
# in C8/D8 move one column to the right, to D8/E8 respectively.
$oldC8 = $ws.Cells.Item(8, 3).Value2
$oldD8 = $ws.Cells.Item(8, 4).Value2

$ws.Cells.Item(8, 5).Value2 = $oldD8
$ws.Cells.Item(8, 4).Value2 = $oldC8

$ws.Cells.Item(8, 3).Value2 = 3.7415500000000002
$ws.Cells.Item(8, 3).NumberFormat = "0.0000"

# Cosmetic: leave the cursor where the author left it when they saved.
$ws.Range("E14").Select()
